$d = $word.ActiveDocument

$replacements = @(
    @{old="65×58=3770"; new="29×28=812"},
    @{old="70×81=5670"; new="30×75=2250"},
    @{old="49×56=2744"; new="58×28=1624"},
    @{old="44×31=1364"; new="24×97=2328"},
    @{old="74×90=6660"; new="15×96=1440"},
    @{old="60×34=2040"; new="56×59=3304"},
    @{old="61×47=2867"; new="52×45=2340"},
    @{old="26×64=1664"; new="62×50=3100"},
    @{old="39×49=1911"; new="47×76=3572"},
    @{old="75×66=4950"; new="62×45=2790"},
    @{old="43×93=3999"; new="12×83=996"},
    @{old="61×24=1464"; new="80×27=2160"},
    @{old="13×39=507"; new="52×22=1144"},
    @{old="71×79=5609"; new="81×11=891"},
    @{old="95×98=9310"; new="45×81=3645"},
    @{old="53×74=3922"; new="44×18=792"},
    @{old="50×29=1450"; new="65×44=2860"},
    @{old="61×90=5490"; new="79×58=4582"},
    @{old="76×61=4636"; new="24×60=1440"},
    @{old="83×84=6972"; new="77×68=5236"},
    @{old="52×49=2548"; new="19×43=817"},
    @{old="18×71=1278"; new="85×66=5610"},
    @{old="87×41=3567"; new="24×56=1344"},
    @{old="41×78=3198"; new="52×22=1144"},
    @{old="45×65=2925"; new="27×35=945"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
